# merge de archivos por mala configuracion
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# STATUS column (F4/F5) currently shows "Cerrada" -> change to "En proceso"
$ws.Range("F4").Value = "En proceso"
$ws.Range("F5").Value = "En proceso"

# Clear "FECHA REAL DE CIERRE" (E4/E5) values, leaving the cells formatted but empty
$ws.Range("E4").ClearContents()
$ws.Range("E5").ClearContents()

# Scroll/selection state: topLeftCell A1 -> C1, selection A5 -> G5
$ws.Range("G5").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
